$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new I0 and IF columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the existing header row (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J77
$data = @{
    2 = @(8, 8)
    3 = @(8, 8)
    4 = @(7, 8)
    5 = @(8, 8)
    6 = @(8, 8)
    7 = @(8, 8)
    8 = @(8, 8)
    9 = @(7, 8)
    10 = @(8, 8)
    11 = @(8, 8)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(8, 8)
    16 = @(8, 8)
    17 = @(8, 8)
    18 = @(7, 8)
    19 = @(8, 8)
    20 = @(8, 8)
    21 = @(8, 8)
    22 = @(8, 8)
    23 = @(8, 8)
    24 = @(8, 8)
    25 = @(8, 8)
    26 = @(8, 8)
    27 = @(9, 9)
    28 = @(8, 8)
    29 = @(8, 8)
    30 = @(8, 8)
    31 = @(8, 8)
    32 = @(8, 8)
    33 = @(7, 8)
    34 = @(8, 8)
    35 = @(8, 8)
    36 = @(11, 11)
    37 = @(8, 8)
    38 = @(8, 8)
    39 = @(8, 8)
    40 = @(8, 8)
    41 = @(8, 8)
    42 = @(8, 8)
    43 = @(9, 9)
    44 = @(8, 8)
    45 = @(8, 8)
    46 = @(7, 7)
    47 = @(8, 8)
    48 = @(8, 8)
    49 = @(8, 8)
    50 = @(7, 7)
    51 = @(8, 8)
    52 = @(7, 7)
    53 = @(8, 8)
    54 = @(7, 7)
    55 = @(8, 8)
    56 = @(8, 8)
    57 = @(9, 9)
    58 = @(8, 8)
    59 = @(8, 8)
    60 = @(8, 8)
    61 = @(8, 8)
    62 = @(7, 7)
    63 = @(8, 9)
    64 = @(8, 8)
    65 = @(7, 7)
    66 = @(8, 8)
    67 = @(6, 7)
    68 = @(7, 7)
    69 = @(6, 7)
    70 = @(6, 6)
    71 = @(5, 5)
    72 = @(4, 4)
    73 = @(7, 7)
    74 = @(5, 6)
    75 = @(7, 7)
    76 = @(8, 8)
    77 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $pair = $data[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
